$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 2: A2="a" B2="aa" C2="desc" D2=01/01/2026 E2=12/31/2027, clear F2
$ws.Range("A2").Value = "a"
$ws.Range("B2").Value = "aa"
$ws.Range("C2").Value = "desc"
$ws.Range("D2").Value = "01/01/2026"
$ws.Range("E2").Value = "12/31/2027"
$ws.Range("F2").ClearContents()

# Row 3: A3="b" B3="bb" C3="desc", clear D3, E3=12/31/2027, F3="pepo"
$ws.Range("A3").Value = "b"
$ws.Range("B3").Value = "bb"
$ws.Range("C3").Value = "desc"
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "12/31/2027"
$ws.Range("F3").Value = "pepo"

# Row 4: clear A4, B4="cc" C4="desc" D4=01/01/2026 E4=12/31/2027 F4="pepe"
$ws.Range("A4").ClearContents()
$ws.Range("B4").Value = "cc"
$ws.Range("C4").Value = "desc"
$ws.Range("D4").Value = "01/01/2026"
$ws.Range("E4").Value = "12/31/2027"
$ws.Range("F4").Value = "pepe"

# Row 5: A5="efe" B5="dd", clear C5, D5=01/01/2026 E5=12/31/2027 F5="pepe"
$ws.Range("A5").Value = "efe"
$ws.Range("B5").Value = "dd"
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = "01/01/2026"
$ws.Range("E5").Value = "12/31/2027"
$ws.Range("F5").Value = "pepe"

# Update the active cell selection to D10
$ws.Range("D10").Select()
